$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename every occurrence of "CHX_Adapted" to "Adapted"
#    (these live in column F / "Genotype") by scanning the used range.
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Text -eq "CHX_Adapted") {
            $cell.Value = "Adapted"
        }
    }
}

# 2. Clear the stray empty formatted cells left over in columns B and F
#    for rows 38-45 (only the G column keeps a cell there).
$ws.Range("B38:B45").ClearContents()
$ws.Range("F38:F45").ClearContents()

# 3. Move the active selection from I20 to J14
[void]$ws.Range("J14").Select()
